# Updates prices/volumes/dates for "Hortaliza, Vega Central Mapocho de Santiago - Achicoria"
# (Fruta / hortaliza, semanal refresh). Column order per row: D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44438
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = 5500
$ws.Range("P2").Value = 344
# Row 3
$ws.Range("D3").Value = 44363
$ws.Range("J3").Value = 160
$ws.Range("K3").Value = 5500
$ws.Range("M3").Value = 5750
$ws.Range("P3").Value = 359
# Row 4
$ws.Range("D4").Value = 44355
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = 6000
$ws.Range("P4").Value = 375
# Row 5
$ws.Range("D5").Value = 44467
$ws.Range("J5").Value = 52
$ws.Range("K5").Value = 5000
$ws.Range("M5").Value = 5500
$ws.Range("P5").Value = 344
# Row 6
$ws.Range("D6").Value = 44455
$ws.Range("J6").Value = 52
$ws.Range("K6").Value = 5000
$ws.Range("M6").Value = 5500
$ws.Range("P6").Value = 344
# Row 7
$ws.Range("D7").Value = 44313
$ws.Range("J7").Value = 34
# Row 8
$ws.Range("D8").Value = 44358
$ws.Range("J8").Value = 52
# Row 9
$ws.Range("D9").Value = 44403
$ws.Range("J9").Value = 43
# Row 10
$ws.Range("D10").Value = 44350
$ws.Range("J10").Value = 25
# Row 11
$ws.Range("D11").Value = 44474
# Row 12
$ws.Range("D12").Value = 44341
$ws.Range("J12").Value = 51
$ws.Range("K12").Value = 5500
$ws.Range("M12").Value = 5755
$ws.Range("P12").Value = 360
# Row 13
$ws.Range("D13").Value = 44308
$ws.Range("J13").Value = 70
$ws.Range("K13").Value = 6000
$ws.Range("M13").Value = 6000
$ws.Range("P13").Value = 375
# Row 14
$ws.Range("D14").Value = 44371
$ws.Range("J14").Value = 34
$ws.Range("K14").Value = 5500
$ws.Range("M14").Value = 5750
$ws.Range("P14").Value = 359
# Row 15
$ws.Range("D15").Value = 44442
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 6000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 6480
$ws.Range("P15").Value = 405
# Row 16
$ws.Range("D16").Value = 44477
# Row 17
$ws.Range("D17").Value = 44306
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 6000
$ws.Range("M17").Value = 6000
$ws.Range("P17").Value = 375
# Row 18
$ws.Range("D18").Value = 44407
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 5500
$ws.Range("M18").Value = 5744
$ws.Range("P18").Value = 359
# Row 20
$ws.Range("D20").Value = 44330
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 6000
$ws.Range("M20").Value = 6000
$ws.Range("P20").Value = 375
